# Arquitectos.xlsx update
# - header "Nit_empresa"/"Cedula_RL" columns removed, "Nivel_educativo" moves to column F
# - columns D/E/F get new explicit widths
# - a new hidden "Hoja2" sheet is added holding the list of education levels
#   (Arquitecto/a, Técnico/a, Tecnólogo/a)
# - column F (from row 2 down) gets a list-type data validation that pulls its
#   source list from Hoja2!$B$2:$B$4

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Hoja1: header row -------------------------------------------------
# F1 used to hold "Nit_empresa"; it now holds "Nivel_educativo".
$ws1.Range("F1").Value = "Nivel_educativo"
# The old G1 ("Nivel_educativo") / H1 ("Cedula_RL") headers are gone.
$ws1.Range("G1:H1").ClearContents() | Out-Null

# --- Hoja1: column widths -----------------------------------------------
$ws1.Columns.Item(4).ColumnWidth = 17.6
$ws1.Columns.Item(5).ColumnWidth = 14.8
$ws1.Columns.Item(6).ColumnWidth = 17.6

# --- Hoja2: new hidden sheet with the education-level catalog ----------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Hoja2"
$ws2.Range("B2").Value = "Arquitecto/a"
$ws2.Range("B3").Value = "Técnico/a"
$ws2.Range("B4").Value = "Tecnólogo/a"
$ws2.Range("B5").Select() | Out-Null
$ws2.Visible = 0

# --- Hoja1: list data validation on F2:F500 -----------------------------
$ws1.Select() | Out-Null
$rng = $ws1.Range("F2:F500")
$rng.Validation.Add(3, 1, 1, "=Hoja2!`$B`$2:`$B`$4") | Out-Null
$rng.Validation.IgnoreBlank = 1
$rng.Validation.InCellDropdown = 1
$rng.Validation.ShowInput = 1
$rng.Validation.ShowError = 1

# Leave the same cell selected/active as in the final workbook.
$ws1.Range("J13").Select() | Out-Null
